$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated / newly-populated cell values per the target diff
$ws.Range("P2").Value = 20000
$ws.Range("S2").Value = 90862
$ws.Range("T2").Value = 5539.9085
$ws.Range("U2").Value = 65000
$ws.Range("P3").Value = 10000
$ws.Range("S3").Value = 71987
$ws.Range("T3").Value = 5343.505999999999
$ws.Range("U3").Value = 55000
$ws.Range("V3").Value = 4465.744597222222
$ws.Range("W3").Value = 4.961304639950756
$ws.Range("B4").Value = 25000
$ws.Range("I4").Value = 20000
$ws.Range("S4").Value = 69101
$ws.Range("T4").Value = 5250.640500000001
$ws.Range("U4").Value = 45000
$ws.Range("B5").Value = 25000
$ws.Range("I5").Value = 20000
$ws.Range("S5").Value = 67174
$ws.Range("T5").Value = 5194.129499999999
$ws.Range("U5").Value = 45000
$ws.Range("S6").Value = 67709
$ws.Range("T6").Value = 5143.397
$ws.Range("S7").Value = 79680
$ws.Range("T7").Value = 5188.8305
$ws.Range("S8").Value = 66217
$ws.Range("T8").Value = 5628.136500000001
$ws.Range("S9").Value = 73755
$ws.Range("T9").Value = 6549.074000000001
$ws.Range("B10").Value = 12500
$ws.Range("I10").Value = 10000
$ws.Range("P10").Value = 0
$ws.Range("S10").Value = 82980
$ws.Range("T10").Value = 7809.641
$ws.Range("U10").Value = 22500
$ws.Range("B11").Value = 12500
$ws.Range("P11").Value = 0
$ws.Range("S11").Value = 95273
$ws.Range("T11").Value = 12977.8705
$ws.Range("U11").Value = 32500
$ws.Range("P12").Value = 10000
$ws.Range("S12").Value = 97940
$ws.Range("T12").Value = 14879.494
$ws.Range("U12").Value = 55000
$ws.Range("S13").Value = 95015
$ws.Range("T13").Value = 15334.284
$ws.Range("S14").Value = 103218
$ws.Range("T14").Value = 15538.4565
$ws.Range("S15").Value = 100719
$ws.Range("T15").Value = 15473.9585
$ws.Range("S16").Value = 99562
$ws.Range("T16").Value = 15459.3285
$ws.Range("S17").Value = 98100
$ws.Range("T17").Value = 15733.5675
$ws.Range("S18").Value = 94912
$ws.Range("T18").Value = 15999.228
$ws.Range("S19").Value = 93211
$ws.Range("T19").Value = 15687.07
$ws.Range("S20").Value = 95142
$ws.Range("T20").Value = 15318.933
$ws.Range("S21").Value = 91234
$ws.Range("T21").Value = 13434.6415
$ws.Range("S22").Value = 88757
$ws.Range("T22").Value = 11609.815
$ws.Range("S23").Value = 85054
$ws.Range("T23").Value = 9103.451000000001
$ws.Range("S24").Value = 80176
$ws.Range("T24").Value = 6290.4415
$ws.Range("S25").Value = 76785
$ws.Range("T25").Value = 5217.8665

# Cells cleared (values removed) per the target diff
$ws.Range("V18").ClearContents()
$ws.Range("W18").ClearContents()
